$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.969.01'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.794.45'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9999'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5315'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3934'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.36%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07460'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.47'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.18%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.085'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.0000'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.185'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.30%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.478'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.792.58'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001059'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06577'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.21'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.01%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.963'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.993.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.40%  '
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '156.75'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.000.79'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.305'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '122.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1086'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.62%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.092'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.672'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.505'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.73%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07071'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.65%  '
$ws.Range('E36').Value = '  -1.98%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.145'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02275'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.398'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('E40').Value = '  -0.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.186'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6120'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.673'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5721'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '125.16'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.181'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.88%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.920'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06804'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.20'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.67%  '
